# ND07.xlsx -> update test fixture:
#  - shared string "ND07" becomes "ND01" (WMT_Extract!C2:C3)
#  - active/selected sheet moves from WMT_Extract to Inst_Reports
#  - WMT_Extract view scrolls/selects to column D / cell AE2
#  - Inst_Reports view selects cell E14
#  - Inst_Reports!G2 gets a value of 2 (fills the blank column gap)

$wb = $excel.ActiveWorkbook

$wsExtract = $wb.Worksheets.Item("WMT_Extract")
$wsInst    = $wb.Worksheets.Item("Inst_Reports")

# --- content: correct the "ND07" label to "ND01" ---
$wsExtract.Range("C2").Value = "ND01"
$wsExtract.Range("C3").Value = "ND01"

# --- content: no more blank column - fill in Inst_Reports G2 ---
$wsInst.Range("G2").Value = 2

# --- view state: WMT_Extract is no longer the active/selected tab ---
$null = $wsExtract.Activate()
$null = $wsExtract.Range("AE2").Select()

# --- view state: Inst_Reports becomes the active/selected tab ---
$null = $wsInst.Activate()
$null = $wsInst.Range("E14").Select()
